$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (rows 2-9)
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

# Clear column D values for rows 2-9 (cells removed entirely)
$ws.Range("D2:D9").ClearContents()

# D13 becomes a plain value instead of a SUBTOTAL formula
$ws.Range("D13").Value = 195

# Update selection
$ws.Range("B3").Select()
